$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, pushing existing rows 60-65 down to 61-66
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new record
$ws.Cells.Item(60, 1).Value = 11
$ws.Cells.Item(60, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(60, 3).Value = "Bíobío"
$ws.Cells.Item(60, 4).Value = 44568
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(61, 4).NumberFormat
$ws.Cells.Item(60, 5).Value = 8
$ws.Cells.Item(60, 6).Value = 100112021
$ws.Cells.Item(60, 7).Value = "Ají"
$ws.Cells.Item(60, 8).Value = "Americana (o)"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 100
$ws.Cells.Item(60, 11).Value = 26000
$ws.Cells.Item(60, 12).Value = 28000
$ws.Cells.Item(60, 13).Value = 27000
$ws.Cells.Item(60, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(60, 15).Value = "Región Metropolitana"
$ws.Cells.Item(60, 16).Value = 1080
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
